# Append two new ticket rows (145 and 146) to the tickets log, mirroring
# the existing rows' layout: columns A-I (Fecha, Hora, WC47 NACP,
# WC48 P5F, WC49 P5H, WV50 FILTER, SPL, Hora de Reparacion,
# Tiempo de Reparacion). The sheet has no date/time number formatting -
# every value (including dates and clock times) is stored as plain text -
# so force a text number format before writing to avoid Excel silently
# reinterpreting "2024-05-21" / "12:50:41" as a date/time serial, then
# drop back to the sheet's default "Normal" style once the literal text
# values are in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A145:I146").NumberFormat = "@"

$ws.Range("A145").Value = "2024-05-21"
$ws.Range("B145").Value = "12:50:41"
$ws.Range("C145").Value = "Fallo tornillo"
$ws.Range("D145").Value = "-"
$ws.Range("E145").Value = "-"
$ws.Range("F145").Value = "-"
$ws.Range("G145").Value = "-"
$ws.Range("H145").Value = "12:50:44"
$ws.Range("I145").Value = "0:00:03"

$ws.Range("A146").Value = "2024-05-21"
$ws.Range("B146").Value = "12:58:21"
$ws.Range("C146").Value = "Etiquetadora21212"
$ws.Range("D146").Value = "-"
$ws.Range("E146").Value = "-"
$ws.Range("F146").Value = "-"
$ws.Range("G146").Value = "-"
$ws.Range("H146").Value = "12:58:26"
$ws.Range("I146").Value = "0:00:05"

$ws.Range("A145:I146").Style = "Normal"
